$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Carrera" column: both rows change from "Negocios" to "Sistemas"
$ws.Range("C2").Value = "Sistemas"
$ws.Range("C3").Value = "Sistemas"

# Column widths for B (Periodo) and C (Carrera)
$ws.Columns.Item(2).ColumnWidth = 13.5
$ws.Columns.Item(3).ColumnWidth = 12.5

# Update the active selection to D5
$ws.Range("D5").Select() | Out-Null
